$wb = $excel.ActiveWorkbook

# --- external sheet: N7 from -70 to -60 ---
$wsExternal = $wb.Worksheets.Item("external")
$wsExternal.Range("N7").Value = -60
$wsExternal.Range("M29").Select()

# --- internal sheet: N7 from 70 to 60 ---
$wsInternal = $wb.Worksheets.Item("internal")
$wsInternal.Range("N7").Value = 60
$wsInternal.Range("L35").Select()

$excel.CalculateFullRebuild()
